# Apply the "updated on 20 Nov" edit to the menstrual cups tracking workbook.

$wb = $excel.ActiveWorkbook

# 1. Rename the first sheet from "Sheet1_2(Nov 19)" to "Sheet1_2(Nov 20)"
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Sheet1_2(Nov 20)"

# Make sure this sheet is the active one (it already is in the source file).
$ws.Activate()

# 2. Fill in the new "Nov 20" readings in column K for rows 8-16.
# K8 already exists (blank, formatted) - just set its value.
$ws.Range("K8").Value = 2

# K9:K16 are brand new cells - copy the formatting from the neighbouring
# "Nov 19" column (J) before putting the value in, so the new cells pick
# up the same style used throughout that column.
$ws.Range("J9").Copy()
$ws.Range("K9").PasteSpecial(-4122)
$ws.Range("K9").Value = 125

$ws.Range("J10").Copy()
$ws.Range("K10").PasteSpecial(-4122)
$ws.Range("K10").Value = 279

$ws.Range("J11").Copy()
$ws.Range("K11").PasteSpecial(-4122)
$ws.Range("K11").Value = 1097

$ws.Range("J12").Copy()
$ws.Range("K12").PasteSpecial(-4122)
$ws.Range("K12").Value = 697

$ws.Range("J13").Copy()
$ws.Range("K13").PasteSpecial(-4122)
$ws.Range("K13").Value = 223

$ws.Range("J14").Copy()
$ws.Range("K14").PasteSpecial(-4122)
$ws.Range("K14").Value = 614

$ws.Range("J15").Copy()
$ws.Range("K15").PasteSpecial(-4122)
$ws.Range("K15").Value = 211

$ws.Range("J16").Copy()
$ws.Range("K16").PasteSpecial(-4122)
$ws.Range("K16").Value = 1176

# 3. Move/update the selection to K16, matching the saved cursor position.
$ws.Range("K16").Select()
